$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G
# (row 2..10 shift up from the row below in the previous version; row 11 is new data)
$data = @{
    2  = @(0.03281260289867554, 0.712211698162994,  1.276354782257036,  1.12975872745336,   1.162021718847977, 18)
    3  = @(-0.08894343397640261, 0.6850734014878143, 1.351404884049018,  1.162499412494053,  1.194764546449047, 17)
    4  = @(0.01098408408208978, 0.5782217553777625, 0.6081029103588239, 0.7798095346678084, 0.805303924577184,  16)
    5  = @(0.1235330352319591,  0.5509845820580629, 0.7140980823158125, 0.8450432428673769, 0.8653061007871674, 15)
    6  = @(0.08484111210645341, 0.6180194413583581, 0.73762096509949,   0.858848627582003,  0.8869099907129777, 14)
    7  = @(0.09014767153583389, 0.6854765477197997, 0.9025650063113468, 0.9500342132320009, 0.9843652458633849, 13)
    8  = @(0.1184762789139503,  0.603487522271732,  0.7334425100585387, 0.856412581679262,  0.8858930048630905, 12)
    9  = @(0.1655637249968011,  0.6554059886689477, 0.9061585844750577, 0.9519236232361595, 0.9831693451202929, 11)
    10 = @(0.1561833212529775,  0.8228064401698413, 1.151759043685242,  1.073200374433983,  1.119208953704769,  10)
    11 = @(0.2159576617971593,  0.7410366868087477, 1.038095822625857,  1.018869875217565,  1.05612043574778,    9)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
